# Update the cryptos list worksheet with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a "Price" (column D) cell while forcing it to be
# stored as text. Excel's COM layer will happily reinterpret plain numeric-
# looking strings (e.g. "528.64", "0.998", "0.0000133") as real numbers when
# assigned via .Value, which would change the cell's type/formatting. The
# original workbook always stores these as plain text, so we briefly force a
# text number format, assign the value, then restore the default "Normal"
# style so no stray formatting is left behind on the cell.
function Set-TextValue($address, $value) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "58.917.41"
$ws.Range("E2").Value = "  +1.28%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.586.93"
$ws.Range("E3").Value = "  -1.00%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
Set-TextValue "D5" "528.64"
$ws.Range("E5").Value = "  +1.63%  "

# Row 6 - Solana
Set-TextValue "D6" "139.05"
$ws.Range("E6").Value = "  -2.54%  "

# Row 7 - USDC
Set-TextValue "D7" "0.998"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.23%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.598.82"
$ws.Range("E9").Value = "  -0.54%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -1.27%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.04%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.29%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +3.13%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "3.044.98"
$ws.Range("E14").Value = "  -0.71%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "58.886.18"
$ws.Range("E15").Value = "  +1.27%  "

# Row 16 - Avalanche
Set-TextValue "D16" "20.48"
$ws.Range("E16").Value = "  +0.74%  "

# Row 17 & 18 - WrappedEther and ShibaInu swapped order
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0000133"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "2.573.09"
$ws.Range("E18").Value = "  -0.85%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "344.26"
$ws.Range("E19").Value = "  +1.65%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  -0.34%  "

# Row 21 - Chainlink
Set-TextValue "D21" "10.06"
$ws.Range("E21").Value = "  -1.31%  "

# Row 22 - Uniswap
Set-TextValue "D22" "6.41"
$ws.Range("E22").Value = "  -0.61%  "

# Row 23 - (stablecoin)
Set-TextValue "D23" "0.998"
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
Set-TextValue "D24" "67.18"
$ws.Range("E24").Value = "  +2.68%  "

# Row 25
Set-TextValue "D25" "0.166"
$ws.Range("E25").Value = "  -0.63%  "

# Row 26
$ws.Range("E26").Value = "  +0.39%  "

# Row 27
Set-TextValue "D27" "0.999"
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("E28").Value = "  +0.26%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("E30").Value = "  -3.06%  "

# Row 31
$ws.Range("E31").Value = "  +1.35%  "

# Row 32
$ws.Range("E32").Value = "  -3.80%  "

# Row 33
Set-TextValue "D33" "18.69"
$ws.Range("E33").Value = "  -0.22%  "

# Row 34 - Monero
Set-TextValue "D34" "149.10"
$ws.Range("E34").Value = "  -0.16%  "

# Row 35
$ws.Range("E35").Value = "  -1.00%  "

# Row 36
$ws.Range("E36").Value = "  -1.44%  "

# Row 37 - OKB
Set-TextValue "D37" "36.68"
$ws.Range("E37").Value = "  +1.68%  "

# Row 38
$ws.Range("E38").Value = "  +0.80%  "

# Row 39
$ws.Range("E39").Value = "  -3.74%  "

# Row 40 - SuiNetwork
Set-TextValue "D40" "0.807"
$ws.Range("E40").Value = "  -4.96%  "

# Row 41
$ws.Range("E41").Value = "  -0.24%  "

# Row 42 - FirstDigitalUSD
Set-TextValue "D42" "0.998"
$ws.Range("E42").Value = "  -0.06%  "

# Row 43 & 44 - Mantle and WhiteBITCoin swapped order
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D43" "10.78"
$ws.Range("E43").Value = "  +0.99%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D44" "0.597"
$ws.Range("E44").Value = "  -1.62%  "

# Row 45 - Bittensor
Set-TextValue "D45" "267.42"
$ws.Range("E45").Value = "  -1.08%  "

# Row 46
$ws.Range("E46").Value = "  -0.10%  "

# Row 47
$ws.Range("E47").Value = "  -1.41%  "

# Row 48 - EnergySwap
Set-TextValue "D48" "18.25"
$ws.Range("E48").Value = "  -2.36%  "

# Row 49 - Maker
Set-TextValue "D49" "1.956.50"
$ws.Range("E49").Value = "  -0.30%  "

# Row 50
$ws.Range("E50").Value = "  -0.35%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "18.11"
$ws.Range("E51").Value = "  -0.36%  "
